$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two "Vostro/Dell" enrollment device rows (ids 589 and 638);
# this shifts the following rows (736.. ) up by two and drops the now
# unused shared strings (Vostro, Dell, DKS, To take enrollments + Arabic
# translations) automatically.
$ws.Rows("6:7").Delete()

# Restore the selection that was saved with the workbook.
$ws.Range("E16").Select()

# Touch page setup so a print-settings relationship/element is written.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
